$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "['BestOffice High-Back Gaming Chair PC Office Chair Computer Racing Chair PU Desk Task Chair Ergonomic Executive Swivel Rolling Chair with Lumbar Support for Back Pain Women, Men,White']"
$ws.Range("B2").Value = "['4.4']"
$ws.Range("C2").Value = "['`$89.99']"
$ws.Range("D2").Value = "['https://m.media-amazon.com/images/I/61t4mpabO+L._AC_SY450_.jpg']"
$ws.Range("E2").Value = "['About this itemPu LeatherCOMFORT FROM EVERY ANGLE - Our high-back gaming chair is thickly cushioned for maximum comfort, whether you’re spending long hours at the office, in front of the computer, or gaming. Covered with breathable premium PU leather, the chair is equipped with freely adjustable lumbar support and a headrest pillow to protect your spine and neck, while the adjustable reclining feature allows users to lock in any reclining position, from 90-135°.MAXIMUM MOBILITY - Our chair are designed with long-lasting ergonomic construction, equipped with a BIFMA-certified heavy-duty metal base, providing stability and mobility with a 360° range of motion on smooth-rolling nylon casters so you can glide easily between work stations.EASY ASSEMBLY – Our chair comes ready to assemble, with all the hardware and necessary tools. With step-by-step instructions, you’ll be set up and ready to game, take on the office in about 10-15 minutes!CUSTOMER GUARANTEE – We want all of our customers to feel 100% satisfied, and ready to unwind from the comfort of our High-Back Office Chair. If you’re not satisfied, or have any questions, please contact BestOffice customer service.BUILT TO LAST - The high-back chair is built to endure, with a heavy-duty metal base, offering great stability and mobility for a capacity of up to 250 lbs.']"
$ws.Range("F2").Value = "['Comfort, Ease of assembly, Value, Appearance, Quality, Support, Size, Reclining']"
